$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(6).Delete()
$ws.Range("F1").Select()
